# Case with 380 kV done: update pl_mw result values for rows 2-25
# (columns B, C, D, F, G, I, L, O) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.338543043387347
$ws.Cells.Item(2, 3).Value = 0.2813589074824279
$ws.Cells.Item(2, 4).Value = 0.02398733850719026
$ws.Cells.Item(2, 6).Value = 0.6040807872797771
$ws.Cells.Item(2, 7).Value = 0.002410267303696152
$ws.Cells.Item(2, 9).Value = 0.5401732232996466
$ws.Cells.Item(2, 12).Value = 0.3007798078620425
$ws.Cells.Item(2, 15).Value = 2.023845372558668
$ws.Cells.Item(3, 2).Value = 1.195971690265196
$ws.Cells.Item(3, 3).Value = 0.2639301899934594
$ws.Cells.Item(3, 4).Value = 0.02211355234758372
$ws.Cells.Item(3, 6).Value = 0.6038970663067218
$ws.Cells.Item(3, 7).Value = 0.002413129957934696
$ws.Cells.Item(3, 9).Value = 0.5503736551894391
$ws.Cells.Item(3, 12).Value = 0.2893174225126529
$ws.Cells.Item(3, 15).Value = 2.037845517922548
$ws.Cells.Item(4, 2).Value = 1.108304077494893
$ws.Cells.Item(4, 3).Value = 0.2531897416896243
$ws.Cells.Item(4, 4).Value = 0.02095752483778313
$ws.Cells.Item(4, 6).Value = 0.6043341583847521
$ws.Cells.Item(4, 7).Value = 0.002414981385356337
$ws.Cells.Item(4, 9).Value = 0.5571156514128255
$ws.Cells.Item(4, 12).Value = 0.2824536761256695
$ws.Cells.Item(4, 15).Value = 2.048296630088643
$ws.Cells.Item(5, 2).Value = 1.072548741422565
$ws.Cells.Item(5, 3).Value = 0.2488034291330905
$ws.Cells.Item(5, 4).Value = 0.02048508053957221
$ws.Cells.Item(5, 6).Value = 0.6046502408421972
$ws.Cells.Item(5, 7).Value = 0.002415759501094551
$ws.Cells.Item(5, 9).Value = 0.5599832255359409
$ws.Cells.Item(5, 12).Value = 0.2797004615133858
$ws.Cells.Item(5, 15).Value = 2.0530208597269
$ws.Cells.Item(6, 2).Value = 1.066609843196773
$ws.Cells.Item(6, 3).Value = 0.2480745223288068
$ws.Cells.Item(6, 4).Value = 0.02040655074555531
$ws.Cells.Item(6, 6).Value = 0.6047110496029475
$ws.Cells.Item(6, 7).Value = 0.002415890136846016
$ws.Cells.Item(6, 9).Value = 0.5604666333829691
$ws.Cells.Item(6, 12).Value = 0.2792459392174607
$ws.Cells.Item(6, 15).Value = 2.053833385271048
$ws.Cells.Item(7, 2).Value = 1.107821987395027
$ws.Cells.Item(7, 3).Value = 0.2531306243133429
$ws.Cells.Item(7, 4).Value = 0.02095115872523223
$ws.Cells.Item(7, 6).Value = 0.604337862988956
$ws.Cells.Item(7, 7).Value = 0.002414991783495257
$ws.Cells.Item(7, 9).Value = 0.5571538383664887
$ws.Cells.Item(7, 12).Value = 0.2824163678756548
$ws.Cells.Item(7, 15).Value = 2.048358460158212
$ws.Cells.Item(8, 2).Value = 1.289412306359395
$ws.Cells.Item(8, 3).Value = 0.2753578143520485
$ws.Cells.Item(8, 4).Value = 0.02334242096022621
$ws.Cells.Item(8, 6).Value = 0.6039031193650075
$ws.Cells.Item(8, 7).Value = 0.00241123493252708
$ws.Cells.Item(8, 9).Value = 0.5435907805346893
$ws.Cells.Item(8, 12).Value = 0.2967914242606184
$ws.Cells.Item(8, 15).Value = 2.028286942682314
$ws.Cells.Item(9, 2).Value = 1.644415665575309
$ws.Cells.Item(9, 3).Value = 0.3186214386605286
$ws.Cells.Item(9, 4).Value = 0.02798669033025192
$ws.Cells.Item(9, 6).Value = 0.6074288944882511
$ws.Cells.Item(9, 7).Value = 0.002404608354502186
$ws.Cells.Item(9, 9).Value = 0.5208050482568005
$ws.Cells.Item(9, 12).Value = 0.3263644371581051
$ws.Cells.Item(9, 15).Value = 2.003697709561408
$ws.Cells.Item(10, 2).Value = 1.904492051990701
$ws.Cells.Item(10, 3).Value = 0.3501952393274337
$ws.Cells.Item(10, 4).Value = 0.03137009517538303
$ws.Cells.Item(10, 6).Value = 0.6127108210455248
$ws.Cells.Item(10, 7).Value = 0.002400186748264673
$ws.Cells.Item(10, 9).Value = 0.5064034727523499
$ws.Cells.Item(10, 12).Value = 0.3489397882751604
$ws.Cells.Item(10, 15).Value = 1.994712014660649
$ws.Cells.Item(11, 2).Value = 2.022631450598681
$ws.Cells.Item(11, 3).Value = 0.3645101100300678
$ws.Cells.Item(11, 4).Value = 0.03290279979650279
$ws.Cells.Item(11, 6).Value = 0.6157030322003791
$ws.Cells.Item(11, 7).Value = 0.002398271326198344
$ws.Cells.Item(11, 9).Value = 0.5003630605875209
$ws.Cells.Item(11, 12).Value = 0.3593952112195495
$ws.Cells.Item(11, 15).Value = 1.992612106758742
$ws.Cells.Item(12, 2).Value = 2.067341470047495
$ws.Cells.Item(12, 3).Value = 0.3699235461115222
$ws.Cells.Item(12, 4).Value = 0.0334822434532569
$ws.Cells.Item(12, 6).Value = 0.6169212269532949
$ws.Cells.Item(12, 7).Value = 0.002397559735745174
$ws.Cells.Item(12, 9).Value = 0.4981494705207012
$ws.Cells.Item(12, 12).Value = 0.3633811587519915
$ws.Cells.Item(12, 15).Value = 1.992104029784088
$ws.Cells.Item(13, 2).Value = 2.057713595480152
$ws.Cells.Item(13, 3).Value = 0.368757995956571
$ws.Cells.Item(13, 4).Value = 0.03335749306386759
$ws.Cells.Item(13, 6).Value = 0.6166550761207645
$ws.Cells.Item(13, 7).Value = 0.002397712379638131
$ws.Cells.Item(13, 9).Value = 0.4986229202895558
$ws.Cells.Item(13, 12).Value = 0.3625215260346835
$ws.Cells.Item(13, 15).Value = 1.992200663525722
$ws.Cells.Item(14, 2).Value = 2.026310320957521
$ws.Cells.Item(14, 3).Value = 0.3649556249330317
$ws.Cells.Item(14, 4).Value = 0.03295049037824072
$ws.Cells.Item(14, 6).Value = 0.6158015461063542
$ws.Cells.Item(14, 7).Value = 0.002398212508337832
$ws.Cells.Item(14, 9).Value = 0.5001794667027717
$ws.Cells.Item(14, 12).Value = 0.3597226023516953
$ws.Cells.Item(14, 15).Value = 1.992564545487483
$ws.Cells.Item(15, 2).Value = 2.007071350641752
$ws.Cells.Item(15, 3).Value = 0.3626256022101018
$ws.Cells.Item(15, 4).Value = 0.0327010636566385
$ws.Cells.Item(15, 6).Value = 0.6152898283978487
$ws.Cells.Item(15, 7).Value = 0.002398520638877458
$ws.Cells.Item(15, 9).Value = 0.5011425142541448
$ws.Cells.Item(15, 12).Value = 0.3580116583825657
$ws.Cells.Item(15, 15).Value = 1.992824862107682
$ws.Cells.Item(16, 2).Value = 1.896767736125753
$ws.Cells.Item(16, 3).Value = 0.3492587276998904
$ws.Cells.Item(16, 4).Value = 0.03126979710041411
$ws.Cells.Item(16, 6).Value = 0.6125271638672487
$ws.Cells.Item(16, 7).Value = 0.002400313851556146
$ws.Cells.Item(16, 9).Value = 0.5068085358336596
$ws.Cells.Item(16, 12).Value = 0.3482602393936958
$ws.Cells.Item(16, 15).Value = 1.994889356280964
$ws.Cells.Item(17, 2).Value = 1.829054712202435
$ws.Cells.Item(17, 3).Value = 0.3410459764474751
$ws.Cells.Item(17, 4).Value = 0.03039009150894856
$ws.Cells.Item(17, 6).Value = 0.6109835765280422
$ws.Cells.Item(17, 7).Value = 0.002401438468042329
$ws.Cells.Item(17, 9).Value = 0.5104155410963997
$ws.Cells.Item(17, 12).Value = 0.3423256361576819
$ws.Cells.Item(17, 15).Value = 1.996665873966464
$ws.Cells.Item(18, 2).Value = 1.790092043017125
$ws.Cells.Item(18, 3).Value = 0.3363177028329289
$ws.Cells.Item(18, 4).Value = 0.02988350558648278
$ws.Cells.Item(18, 6).Value = 0.6101512057150629
$ws.Cells.Item(18, 7).Value = 0.002402094356641375
$ws.Cells.Item(18, 9).Value = 0.5125382683631017
$ws.Cells.Item(18, 12).Value = 0.3389296966198287
$ws.Cells.Item(18, 15).Value = 1.997874698331145
$ws.Cells.Item(19, 2).Value = 1.776897286387793
$ws.Cells.Item(19, 3).Value = 0.3347160260553323
$ws.Cells.Item(19, 4).Value = 0.02971188187419216
$ws.Cells.Item(19, 6).Value = 0.609878893862458
$ws.Cells.Item(19, 7).Value = 0.002402317983828573
$ws.Cells.Item(19, 9).Value = 0.5132652344881699
$ws.Cells.Item(19, 12).Value = 0.3377828935943938
$ws.Cells.Item(19, 15).Value = 1.998316067012041
$ws.Cells.Item(20, 2).Value = 1.836264544259905
$ws.Cells.Item(20, 3).Value = 0.3419207078971453
$ws.Cells.Item(20, 4).Value = 0.03048380028317865
$ws.Cells.Item(20, 6).Value = 0.6111421516175355
$ws.Cells.Item(20, 7).Value = 0.002401317815914716
$ws.Cells.Item(20, 9).Value = 0.510026591572899
$ws.Cells.Item(20, 12).Value = 0.3429555750071955
$ws.Cells.Item(20, 15).Value = 1.996457395634849
$ws.Cells.Item(21, 2).Value = 2.035534969390653
$ws.Cells.Item(21, 3).Value = 0.3660726743197245
$ws.Cells.Item(21, 4).Value = 0.03307006321631434
$ws.Cells.Item(21, 6).Value = 0.6160499359472595
$ws.Cells.Item(21, 7).Value = 0.002398065236067938
$ws.Cells.Item(21, 9).Value = 0.4997202663684028
$ws.Cells.Item(21, 12).Value = 0.3605439893715072
$ws.Cells.Item(21, 15).Value = 1.992449862211885
$ws.Cells.Item(22, 2).Value = 2.165612425502786
$ws.Cells.Item(22, 3).Value = 0.3818146956950557
$ws.Cells.Item(22, 4).Value = 0.03475473389321593
$ws.Cells.Item(22, 6).Value = 0.6197536994883706
$ws.Cells.Item(22, 7).Value = 0.002396019534236907
$ws.Cells.Item(22, 9).Value = 0.4934146942099957
$ws.Cells.Item(22, 12).Value = 0.372194692444026
$ws.Cells.Item(22, 15).Value = 1.99150480921486
$ws.Cells.Item(23, 2).Value = 2.096202773732387
$ws.Cells.Item(23, 3).Value = 0.3734169069597897
$ws.Cells.Item(23, 4).Value = 0.03385611736348437
$ws.Cells.Item(23, 6).Value = 0.6177314101501992
$ws.Cells.Item(23, 7).Value = 0.002397104060774567
$ws.Cells.Item(23, 9).Value = 0.4967406302247319
$ws.Cells.Item(23, 12).Value = 0.3659622533231186
$ws.Cells.Item(23, 15).Value = 1.991855600020358
$ws.Cells.Item(24, 2).Value = 1.833005083522778
$ws.Cells.Item(24, 3).Value = 0.3415252625074459
$ws.Cells.Item(24, 4).Value = 0.03044143713383107
$ws.Cells.Item(24, 6).Value = 0.6110702883428587
$ws.Cells.Item(24, 7).Value = 0.002401372333657596
$ws.Cells.Item(24, 9).Value = 0.510202283004709
$ws.Cells.Item(24, 12).Value = 0.3426707299947367
$ws.Cells.Item(24, 15).Value = 1.996551064788775
$ws.Cells.Item(25, 2).Value = 1.548503122685588
$ws.Cells.Item(25, 3).Value = 0.306953671795668
$ws.Cells.Item(25, 4).Value = 0.02673524591683218
$ws.Cells.Item(25, 6).Value = 0.60600402732954
$ws.Cells.Item(25, 7).Value = 0.002406322205244317
$ws.Cells.Item(25, 9).Value = 0.5265595451710361
$ws.Cells.Item(25, 12).Value = 0.3182155594003859
$ws.Cells.Item(25, 15).Value = 2.008760619824443
